$d = $word.ActiveDocument

# --- helpers -------------------------------------------------------------

function Wrap-Xml($inner) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $inner +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $xml
}

# Replace paragraph $index's content (keeping it as a single paragraph) with
# the raw w:p inner-xml supplied in $innerParaXml (a full "<w:p>...</w:p>").
function Set-ParaXml($index, $innerParaXml) {
    $xml = Wrap-Xml $innerParaXml
    $d.Paragraphs($index).Range.InsertXML($xml)
}

function Find-ParaIndex($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

$HL = '<w:rPr><w:highlight w:val="yellow"/></w:rPr>'
$PPR_HL = '<w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>'

# --- District Summary block ----------------------------------------------

$districtIdx = Find-ParaIndex "### District Summary"
if ($districtIdx -eq -1) {
    throw "Could not locate '### District Summary' paragraph"
}

# ### District Summary  (pPr highlight + run highlight)
$idx = $districtIdx
$xml = '<w:p>' + $PPR_HL + '<w:r>' + $HL + '<w:t>### District Summary</w:t></w:r></w:p>'
Set-ParaXml $idx $xml

# blank paragraph right after (pPr highlight only, no run)
$idx = $districtIdx + 1
$xml = '<w:p>' + $PPR_HL + '</w:p>'
Set-ParaXml $idx $xml

# * Create a high level snapshot ... (runs highlighted, no pPr highlight)
$idx = $districtIdx + 2
$xml = '<w:p>' +
    '<w:r>' + $HL + '<w:t xml:space="preserve">* Create a </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $HL + '<w:t>high level</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $HL + '<w:t xml:space="preserve"> snapshot (in table form) of the district''s key metrics, including:</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * Total Schools  (pPr highlight; split runs "  " + "* Total Schools")
$idx = $districtIdx + 3
$xml = '<w:p>' + $PPR_HL +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r>' + $HL + '<w:t>* Total Schools</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * Total Students  (pPr highlight; single run highlighted)
$idx = $districtIdx + 4
$xml = '<w:p>' + $PPR_HL +
    '<w:r>' + $HL + '<w:t xml:space="preserve">  * Total Students</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * Total Budget  (no pPr highlight; single run highlighted)
$idx = $districtIdx + 5
$xml = '<w:p>' +
    '<w:r>' + $HL + '<w:t xml:space="preserve">  * Total Budget</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * Average Math Score  (pPr highlight; split runs "  " + "* Average Math Score")
$idx = $districtIdx + 6
$xml = '<w:p>' + $PPR_HL +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r>' + $HL + '<w:t>* Average Math Score</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * Average Reading Score  (no pPr highlight; single run highlighted)
$idx = $districtIdx + 7
$xml = '<w:p>' +
    '<w:r>' + $HL + '<w:t xml:space="preserve">  * Average Reading Score</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * % Passing Math ...  (no pPr highlight; split runs "  " + rest)
$idx = $districtIdx + 8
$xml = '<w:p>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r>' + $HL + '<w:t>* % Passing Math (The percentage of students that passed math.)</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * % Passing Reading ...  (no pPr highlight; split runs "  " + rest)
$idx = $districtIdx + 9
$xml = '<w:p>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r>' + $HL + '<w:t>* % Passing Reading (The percentage of students that passed reading.)</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * % Overall Passing ...  (no pPr highlight; split runs "  " + rest)
$idx = $districtIdx + 10
$xml = '<w:p>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r>' + $HL + '<w:t>* % Overall Passing (The percentage of students that passed math **and** reading.)</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

# --- School Summary block -------------------------------------------------

$schoolIdx = Find-ParaIndex "### School Summary"
if ($schoolIdx -eq -1) {
    throw "Could not locate '### School Summary' paragraph"
}

#   * School Name  (pPr highlight; split runs "  " + "* School Name")
$idx = $schoolIdx + 3
$xml = '<w:p>' + $PPR_HL +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r>' + $HL + '<w:t>* School Name</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * School Type  (pPr highlight; single run highlighted)
$idx = $schoolIdx + 4
$xml = '<w:p>' + $PPR_HL +
    '<w:r>' + $HL + '<w:t xml:space="preserve">  * School Type</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

#   * Total Students  (no pPr highlight; single run highlighted)
$idx = $schoolIdx + 5
$xml = '<w:p>' +
    '<w:r>' + $HL + '<w:t xml:space="preserve">  * Total Students</w:t></w:r>' +
    '</w:p>'
Set-ParaXml $idx $xml

Write-Output "edit applied"
